$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert the new row 14 (copies formatting from row 13 above it, which we
# then selectively strip so only the columns that should carry the
# "customFormat" style actually keep it).
$ws.Rows(14).Insert()
$ws.Range("A14:M14").ClearFormats()
$ws.Range("P14").ClearFormats()
$ws.Range("P14").ClearContents()

$ws.Range("A14").Value = "Bernhardt_2018_PROFTHROSOB.SC"
$ws.Range("B14").Value = "Bernhardt, Joey R. and Sunday, Jennifer M. and Thompson, Patrick L. and O'Connor, Mary I."
$ws.Range("C14").Value = "Nonlinear averaging of thermal experience predicts population growth rates in a thermally variable environment"
$ws.Range("D14").Value = "PROCEEDINGS OF THE ROYAL SOCIETY B-BIOLOGICAL SCIENCES"
$ws.Range("E14").Value = "10.1098/rspb.2018.1076"
$ws.Range("F14").Value = 2018
$ws.Range("G14").Value = "As thermal regimes change worldwide, projections of future population and species persistence often require estimates of how population growth rates depend on temperature. These projections rarely account for how temporal variation in temperature can systematically modify growth rates relative to projections based on constant temperatures. Here, we tested the hypothesis that time-averaged population growth rates in fluctuating thermal environments differ from growth rates in constant conditions as a consequence of Jensen's inequality, and that the thermal performance curves (TPCs) describing population growth in fluctuating environments can be predicted quantitatively based on TPCs generated in constant laboratory conditions. With experimental populations of the green alga Tetraselmis tetrahele, we show that nonlinear averaging techniques accurately predicted increased as well as decreased population growth rates in fluctuating thermal regimes relative to constant thermal regimes. We extrapolate from these results to project critical temperatures for population growth and persistence of 89 phytoplankton species in naturally variable thermal environments. These results advance our ability to predict population dynamics in the context of global change."
$ws.Range("H14").Value = "0962-8452"
$ws.Range("I14").Value = "Nonlinear Averaging of Thermal Experience Predicts Population Growth Rates in a Thermally Variable Environment."
$ws.Range("J14").Value = "selected"
$ws.Range("K14").Value = "background"
$ws.Range("L14").Value = "selected"
$ws.Range("M14").Value = "y"
$ws.Range("N14").Value = "y "
$ws.Range("O14").Value = "figure 2"
$ws.Range("Q14").Value = "check supplementary info, constant vs flux "
$ws.Range("S14").Value = "y"

# New (empty) row 15 - only Q15 carries the inherited fill-applied style.
$ws.Rows(15).Insert()
$ws.Range("N15").ClearFormats()
$ws.Range("N15").ClearContents()
$ws.Range("O15").ClearFormats()
$ws.Range("O15").ClearContents()
$ws.Range("S15").ClearFormats()
$ws.Range("S15").ClearContents()

# Sheet view: move the active selection.
$null = $ws.Range("Q21").Select()
